# Generate Report for Handoff
# Updates the localization-status report: the "b.md" file has now been
# handed off (zh-cn + de-de), so its status flips from "Handed back: in
# sync with en-US" to "Ready for handoff" on the Overview sheet and on
# both language sheets; the zh-cn/de-de sheets also get refreshed handoff
# file names/timestamps and an explanatory error detail, and column P on
# every language sheet is widened to fit the new text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is "b.md"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-01 02:46:17"

# ---------------------------------------------------------------------
# zh-cn sheet: row 2 is "a.md", row 3 is "b.md"
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-01 02:46:13"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a30c4917a8e921fe273bc4c1f35f10e61a64c322/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55d97ab1e0d7ec43a026637ac7fa87d227b14820/e2e/b.md."

# widen column P (Error Detail) to fit the new message
$zhcn.Columns.Item(16).ColumnWidth = 39.1476

# ---------------------------------------------------------------------
# de-de sheet: row 2 is "a.md", row 3 is "b.md"
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-01 02:46:17"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a30c4917a8e921fe273bc4c1f35f10e61a64c322/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55d97ab1e0d7ec43a026637ac7fa87d227b14820/e2e/b.md."

# widen column P (Error Detail) to fit the new message
$dede.Columns.Item(16).ColumnWidth = 39.1476
